# "timings for the sync feature"
#
# Sheet1 held two independent open/close timing tables side by side
# (A:C = seconds+frames -> ms, E:G = seconds+frames -> ms). This replaces
# that with a single simplified sync table: raw millisecond "Open" values
# in column A and a computed "Close" (Open + 600ms) in column B, for 5
# events instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused "frames"/"close-sec"/"ms" columns (C, and E:G)
# entirely -- only the A/B open-timing pair survives.
$ws.Range("C1:G7").EntireColumn.Delete()

# Drop the 6th data row (old row 7) -- the new table only has 5 events.
$ws.Rows("7:7").Delete()

# Reset the header row's own row-level formatting, then re-apply bold
# just to the two surviving header cells.
$ws.Rows("1:1").ClearFormats()
$ws.Range("A1:B1").Font.Bold = $true

# New header text.
$ws.Range("A1").Value = "Open"
$ws.Range("B1").Value = "Close"

# New timing data (ms) -- Close is always Open + 600ms.
$ws.Range("A2").Value = 467
$ws.Range("B2").Formula = "=A2+600"

$ws.Range("A3").Value = 2102
$ws.Range("B3").Formula = "=A3+600"

$ws.Range("A4").Value = 3470
$ws.Range("B4").Formula = "=A4+600"

$ws.Range("A5").Value = 5038
$ws.Range("B5").Formula = "=A5+600"

$ws.Range("A6").Value = 6673
$ws.Range("B6").Formula = "=A6+600"

# The old Print_Area named range ($A$1:$G$7) no longer matches the
# shrunken sheet; point it at a broken reference rather than leaving a
# stale range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!Print_Area") {
        $n.RefersTo = "=Sheet1!#REF!"
    }
}

# Leave the cursor where the author last left it.
[void]$ws.Range("B9").Select()
